$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value2 = "0-A-10"
$ws.Cells.Item(1, 2).Value2 = 0
$ws.Cells.Item(1, 3).Value2 = 139.3840749090528
$ws.Cells.Item(1, 4).Value2 = 1
$ws.Cells.Item(1, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(2, 1).Value2 = "0-A-30"
$ws.Cells.Item(2, 2).Value2 = 0
$ws.Cells.Item(2, 3).Value2 = 414.5390494141789
$ws.Cells.Item(2, 4).Value2 = 1
$ws.Cells.Item(2, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(3, 1).Value2 = "0-B-10"
$ws.Cells.Item(3, 2).Value2 = 0
$ws.Cells.Item(3, 3).Value2 = 125.2643488489878
$ws.Cells.Item(3, 4).Value2 = 0.6
$ws.Cells.Item(3, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(4, 1).Value2 = "0-B-30"
$ws.Cells.Item(4, 2).Value2 = 0
$ws.Cells.Item(4, 3).Value2 = 367.0061185843164
$ws.Cells.Item(4, 4).Value2 = 0.6
$ws.Cells.Item(4, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(5, 1).Value2 = "0.5-A-10"
$ws.Cells.Item(5, 2).Value2 = 0
$ws.Cells.Item(5, 3).Value2 = 159.5855541711678
$ws.Cells.Item(5, 4).Value2 = 1
$ws.Cells.Item(5, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(6, 1).Value2 = "0.5-A-30"
$ws.Cells.Item(6, 2).Value2 = 0
$ws.Cells.Item(6, 3).Value2 = 475.8317777903611
$ws.Cells.Item(6, 4).Value2 = 1
$ws.Cells.Item(6, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(7, 1).Value2 = "0.5-B-10"
$ws.Cells.Item(7, 2).Value2 = 0
$ws.Cells.Item(7, 3).Value2 = 144.6409437369407
$ws.Cells.Item(7, 4).Value2 = 0.6
$ws.Cells.Item(7, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(8, 1).Value2 = "0.5-B-30"
$ws.Cells.Item(8, 2).Value2 = 0
$ws.Cells.Item(8, 3).Value2 = 424.4556792034342
$ws.Cells.Item(8, 4).Value2 = 0.6
$ws.Cells.Item(8, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(9, 1).Value2 = "1-A-10"
$ws.Cells.Item(9, 2).Value2 = 0
$ws.Cells.Item(9, 3).Value2 = 180.290437986355
$ws.Cells.Item(9, 4).Value2 = 1
$ws.Cells.Item(9, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(10, 1).Value2 = "1-A-30"
$ws.Cells.Item(10, 2).Value2 = 0
$ws.Cells.Item(10, 3).Value2 = 537.3438803121122
$ws.Cells.Item(10, 4).Value2 = 1
$ws.Cells.Item(10, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(11, 1).Value2 = "1-B-10"
$ws.Cells.Item(11, 2).Value2 = 0
$ws.Cells.Item(11, 3).Value2 = 165.0778722759138
$ws.Cells.Item(11, 4).Value2 = 0.6
$ws.Cells.Item(11, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(12, 1).Value2 = "1-B-30"
$ws.Cells.Item(12, 2).Value2 = 0
$ws.Cells.Item(12, 3).Value2 = 484.4947180537995
$ws.Cells.Item(12, 4).Value2 = 0.6
$ws.Cells.Item(12, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(13, 1).Value2 = "1.1-A-10"
$ws.Cells.Item(13, 2).Value2 = 0
$ws.Cells.Item(13, 3).Value2 = 184.4531888715271
$ws.Cells.Item(13, 4).Value2 = 1
$ws.Cells.Item(13, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(14, 1).Value2 = "1.1-A-30"
$ws.Cells.Item(14, 2).Value2 = 0
$ws.Cells.Item(14, 3).Value2 = 549.6470835210055
$ws.Cells.Item(14, 4).Value2 = 1
$ws.Cells.Item(14, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(15, 1).Value2 = "1.1-B-10"
$ws.Cells.Item(15, 2).Value2 = 0
$ws.Cells.Item(15, 3).Value2 = 169.1802328329443
$ws.Cells.Item(15, 4).Value2 = 0.6
$ws.Cells.Item(15, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(16, 1).Value2 = "1.1-B-30"
$ws.Cells.Item(16, 2).Value2 = 0
$ws.Cells.Item(16, 3).Value2 = 496.5355733394323
$ws.Cells.Item(16, 4).Value2 = 0.6
$ws.Cells.Item(16, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(17, 1).Value2 = "3-A-10"
$ws.Cells.Item(17, 2).Value2 = 0
$ws.Cells.Item(17, 3).Value2 = 264.1581374446722
$ws.Cells.Item(17, 4).Value2 = 1
$ws.Cells.Item(17, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(18, 1).Value2 = "3-A-30"
$ws.Cells.Item(18, 2).Value2 = 0
$ws.Cells.Item(18, 3).Value2 = 785.4795337091788
$ws.Cells.Item(18, 4).Value2 = 1
$ws.Cells.Item(18, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(19, 1).Value2 = "3-B-10"
$ws.Cells.Item(19, 2).Value2 = 1
$ws.Cells.Item(19, 3).Value2 = 247.8959334482108
$ws.Cells.Item(19, 4).Value2 = 0.6
$ws.Cells.Item(19, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(20, 1).Value2 = "3-B-30"
$ws.Cells.Item(20, 2).Value2 = 1
$ws.Cells.Item(20, 3).Value2 = 726.1851268378294
$ws.Cells.Item(20, 4).Value2 = 0.6
$ws.Cells.Item(20, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(21, 1).Value2 = "5-A-10"
$ws.Cells.Item(21, 2).Value2 = 0
$ws.Cells.Item(21, 3).Value2 = 348.5198581214263
$ws.Cells.Item(21, 4).Value2 = 1
$ws.Cells.Item(21, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(22, 1).Value2 = "5-A-30"
$ws.Cells.Item(22, 2).Value2 = 0
$ws.Cells.Item(22, 3).Value2 = 1034.84863733484
$ws.Cells.Item(22, 4).Value2 = 1
$ws.Cells.Item(22, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(23, 1).Value2 = "5-B-10"
$ws.Cells.Item(23, 2).Value2 = 1
$ws.Cells.Item(23, 3).Value2 = 330.8843809593009
$ws.Cells.Item(23, 4).Value2 = 0.6
$ws.Cells.Item(23, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(24, 1).Value2 = "5-B-30"
$ws.Cells.Item(24, 2).Value2 = 3
$ws.Cells.Item(24, 3).Value2 = 968.7254332588374
$ws.Cells.Item(24, 4).Value2 = 0.6
$ws.Cells.Item(24, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(25, 1).Value2 = "6-A-10"
$ws.Cells.Item(25, 2).Value2 = 2
$ws.Cells.Item(25, 3).Value2 = 390.7805850219103
$ws.Cells.Item(25, 4).Value2 = 1
$ws.Cells.Item(25, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(26, 1).Value2 = "6-A-30"
$ws.Cells.Item(26, 2).Value2 = 2
$ws.Cells.Item(26, 3).Value2 = 1160.093012485466
$ws.Cells.Item(26, 4).Value2 = 1
$ws.Cells.Item(26, 5).Value2 = 0.01378870056105845

$ws.Cells.Item(27, 1).Value2 = "6-B-10"
$ws.Cells.Item(27, 2).Value2 = 3
$ws.Cells.Item(27, 3).Value2 = 372.4052109721124
$ws.Cells.Item(27, 4).Value2 = 0.6
$ws.Cells.Item(27, 5).Value2 = 0.006811749575498305

$ws.Cells.Item(28, 1).Value2 = "6-B-30"
$ws.Cells.Item(28, 2).Value2 = 4
$ws.Cells.Item(28, 3).Value2 = 1090.286427512467
$ws.Cells.Item(28, 4).Value2 = 0.6
$ws.Cells.Item(28, 5).Value2 = 0.006811749575498305

# Remove now-unused rows 29-36 (table shrank from 36 to 28 rows)
$ws.Rows("29:36").Delete()
